# Syllable analysis, improvement on spell.txt
#
# This script reshapes the little "Attribute" side-table (which lived in
# columns K:N) one column to the right (L:O), frees up column K, adds a
# per-word "Total" column (J) that sums the Beginning/Penultimate/Last
# counts, adds two new standalone labels (CV / VC) below the attribute
# table, and appends a grand-total row (28) under the main letter table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Move the "Attribute" summary table from K1:N5 to L1:O5, and add the
#    new "Total" header in J1.
# ---------------------------------------------------------------------

$ws.Range("J1").Value = "Total"

$ws.Range("L1").Value = "Attribute"
$ws.Range("M1").Value = "Beginning"
$ws.Range("N1").Value = "Penultimate"
$ws.Range("O1").Value = "Last"

$ws.Range("L2").Value = "Begins with Consonant"
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = 13

$ws.Range("L3").Value = "Begins with Vowel"
$ws.Range("M3").Value = 5
$ws.Range("N3").Value = 3
$ws.Range("O3").Value = 0

$ws.Range("L4").Value = "Double Consonants"
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 2

$ws.Range("L5").Value = "Double Vowels"
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 4

# New standalone labels just below the attribute table.
$ws.Range("L6").Value = "CV"
$ws.Range("L7").Value = "VC"

# Clear out the old column K (the table now starts at L).
$ws.Range("K1:K5").Clear()

# ---------------------------------------------------------------------
# 2. New column J: per-row total of Beginning + Penultimate + Last
#    (columns G:I), rows 2-27. Enter J2 on its own, then fill J3:J27 as a
#    block so Excel treats it as a shared formula, matching a
#    select-and-fill-down workflow.
# ---------------------------------------------------------------------

$ws.Range("J2").Formula = "=SUM(G2:I2)"
$ws.Range("J3:J27").Formula = "=SUM(G3:I3)"

# Bold styling (same look as the other header/total cells, style index 1
# i.e. the "Normal + bold font" cell style) for the Total column and its
# header.
$ws.Range("J1:J27").Font.Bold = $true

# ---------------------------------------------------------------------
# 3. New grand-total row 28 under the letter-frequency table.
# ---------------------------------------------------------------------

$ws.Range("F28").Value = "Total:"
$ws.Range("G28").Formula = "=SUM(G2:G27)"
$ws.Range("H28:I28").Formula = "=SUM(H2:H27)"
$ws.Range("J28").Formula = "=SUM(J2:J27)"

$ws.Range("F28:J28").Font.Bold = $true

# ---------------------------------------------------------------------
# 4. View state: scroll so column G is left-most visible, select M6.
# ---------------------------------------------------------------------

$ws.Range("M6").Select()
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 4
